$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row, DAMSLTag (col I / 9), DialogAct (col J / 10)
$updates = @(
    ,@(10, "sd", "Statement-non-opinion")
    ,@(11, "sv", "Statement-opinion")
    ,@(14, "sd", "Statement-non-opinion")
    ,@(34, "sd", "Statement-non-opinion")
    ,@(40, "sd", "Statement-non-opinion")
    ,@(48, "sv", "Statement-opinion")
    ,@(60, "sd", "Statement-non-opinion")
    ,@(64, "ba", "Appreciation")
    ,@(70, "sv", "Statement-opinion")
    ,@(89, "sv", "Statement-opinion")
    ,@(111, "sv", "Statement-opinion")
    ,@(113, "b", "Acknowledge (Backchannel)")
    ,@(114, "sd", "Statement-non-opinion")
    ,@(124, "sd", "Statement-non-opinion")
    ,@(142, "b", "Acknowledge (Backchannel)")
    ,@(148, "%", "Uninterpretable")
    ,@(156, "sv", "Statement-opinion")
    ,@(164, "sv", "Statement-opinion")
    ,@(169, "sd", "Statement-non-opinion")
    ,@(203, "sd", "Statement-non-opinion")
    ,@(210, "qy", "Yes-No-Question")
    ,@(232, "sd", "Statement-non-opinion")
    ,@(248, "aa", "Agree/Accept")
    ,@(289, "sd", "Statement-non-opinion")
    ,@(293, "sd", "Statement-non-opinion")
    ,@(317, "ba", "Appreciation")
    ,@(318, "b", "Acknowledge (Backchannel)")
    ,@(331, "sd", "Statement-non-opinion")
    ,@(332, "%", "Uninterpretable")
    ,@(348, "sv", "Statement-opinion")
    ,@(362, "sd", "Statement-non-opinion")
    ,@(382, "aa", "Agree/Accept")
    ,@(395, "aa", "Agree/Accept")
    ,@(434, "sv", "Statement-opinion")
    ,@(442, "aa", "Agree/Accept")
    ,@(448, "sd", "Statement-non-opinion")
    ,@(453, "sv", "Statement-opinion")
    ,@(458, "sd", "Statement-non-opinion")
    ,@(472, "sd", "Statement-non-opinion")
    ,@(473, "sd", "Statement-non-opinion")
    ,@(489, "aa", "Agree/Accept")
    ,@(507, "sd", "Statement-non-opinion")
)

foreach ($u in $updates) {
    $r = $u[0]
    $ws.Cells.Item($r, 9).Value = $u[1]
    $ws.Cells.Item($r, 10).Value = $u[2]
}
